$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The table currently ends with the "trial 12" row (row index 13, 1-based,
# since row 1 is the header). The commit adds a new "trial 13" data row by:
#   1) duplicating the old trial-12 row in place (a fresh blank row is
#      inserted immediately before the existing last row, then filled with
#      the old trial-12 values), and
#   2) turning the former last row into the new trial-13 row with its own
#      updated statistics.

$lastRowIndex = $t.Rows.Count
$oldLastRow = $t.Rows.Item($lastRowIndex)

# Insert a new blank row before the current last row, copying its formatting.
$newRow = $t.Rows.Add($oldLastRow)

# Fill the newly inserted row with the (duplicated) old trial-12 data.
$t.Cell($lastRowIndex, 1).Range.Text = "12"
$t.Cell($lastRowIndex, 2).Range.Text = "1"
$t.Cell($lastRowIndex, 3).Range.Text = "0.05"
$t.Cell($lastRowIndex, 4).Range.Text = "7,134.000"
$t.Cell($lastRowIndex, 5).Range.Text = "4,425"
$t.Cell($lastRowIndex, 6).Range.Text = "<0.001"

# The original last row has now shifted down by one; update it in place to
# hold the new trial-13 data.
$newLastRowIndex = $t.Rows.Count
$t.Cell($newLastRowIndex, 1).Range.Text = "13"
$t.Cell($newLastRowIndex, 2).Range.Text = "0.0038"
$t.Cell($newLastRowIndex, 4).Range.Text = "14,760.000"
$t.Cell($newLastRowIndex, 5).Range.Text = "14,640"
$t.Cell($newLastRowIndex, 6).Range.Text = "0.241"

Write-Output ("Final row count: " + $t.Rows.Count)
